# Generate Report for Handoff
# Update status text and timestamps on all sheets, and widen the
# "Status" columns so the new, longer status text fits.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update the status text from "In Translation" to "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Update the handoff timestamps
$wsOverview.Range("G2").Value = "2016-09-04 04:43:51"
$wsDeDe.Range("H2").Value = "2016-09-04 04:43:51"
$wsZhCn.Range("H2").Value = "2016-09-04 04:43:47"

# Widen the status columns to fit the new, longer text (matches the
# width Excel's own "best fit" would produce for the new status string)
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
